# Rename the worksheet from "Sheet1" to "Tasks"
# (Commit message: "migrated to sheetName from sheet index while reading/writing excel")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Tasks"
